$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'320.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.67%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'49.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'14.21%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.259"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'3.54%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07957"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.44%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.575"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'1.420"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'35.89%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.640"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.06%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1299"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'3.59%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1968"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'5.91%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09507"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.92%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04617"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'10.83%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'-0.28%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001324"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'3.72%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.04158"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.28%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005926"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'3.08%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'-0.06%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.434"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.39%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3461"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'3.18%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.228"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-4.88%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'1.69%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.3091"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'0.001314"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.49%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004256"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-4.91%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'0.03%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003535"
$ws.Range("D26").Style = "Normal"
$ws.Range("D38").Value = "'0.02674"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'8.87%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05879"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'11.32%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.01091"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'84.21%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.008001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.14%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1439"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'7.02%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007715"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'4.63%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008683"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'14.88%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3192"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.96%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006625"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.06%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.05499"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'23.67%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003996"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-4.88%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'0.03%"
$ws.Range("E51").Style = "Normal"
